$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 through 11 (the review rows), keeping row 1 (B1 = 0) intact.
$ws.Range("A2:B11").EntireRow.Delete()
